# grades_h19.xlsx revision:
#   - Column B ("grade") and column C ("name") are swapped so that the
#     "name" column (Traditional/Online) becomes "course_type" in column B,
#     and "grade" moves to column C.
#   - Shared strings end up with "name" removed and "course_type" appended.
#   - Column B gets its own (wider) column-width entry, split out from the
#     combined A:C width block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the B/C data columns for every data row (2-208) ---------------
$dataRange = $ws.Range("B2:C208")
$vals = $dataRange.Value()
$rowCount = $vals.GetLength(0)
for ($i = 1; $i -le $rowCount; $i++) {
    $bVal = $vals[$i, 1]
    $cVal = $vals[$i, 2]
    $vals[$i, 1] = $cVal
    $vals[$i, 2] = $bVal
}
$dataRange.Value = $vals

# --- Fix up the header row: B1 "grade" -> "course_type"; C1 "name" -> "grade"
$ws.Range("B1").Value = "course_type"
$ws.Range("C1").Value = "grade"

# --- Give column B its own width (previously merged with A and C) -------
$ws.Columns.Item(2).ColumnWidth = 10.25
